$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Widen column A to fit the new longer labels
$ws.Columns.Item(1).ColumnWidth = 52.85546875

# Set values in the same order the shared-string table was populated:
# A2, A3, A1, B2, B3, then B1 (which reuses the pre-existing "BoilerMakeBoard_I" string).
$ws.Cells.Item(2, 1).Value = "Assembly Part Number/Revision:"
$ws.Cells.Item(3, 1).Value = "Customer / Company Name: "
$ws.Cells.Item(1, 1).Value = "BOM: "
$ws.Cells.Item(2, 2).Value = "Prototype for production"
$ws.Cells.Item(3, 2).Value = "Thomas Kilbride / Purdue University BoilerMake Hackathon"
$ws.Cells.Item(1, 2).Value = "BoilerMakeBoard_I"

# Row heights for the new annotation rows
$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(3).RowHeight = 15.75

# Style the three label cells: bold Arial 12, left/center aligned
$labelRange = $ws.Range("A1:A3")
$labelRange.Style = "Normal 2"
$labelRange.Font.Bold = $true
$labelRange.Font.Name = "Arial"
$labelRange.Font.Size = 12
$labelRange.HorizontalAlignment = -4131  # xlLeft
$labelRange.VerticalAlignment = -4108    # xlCenter

$ws.Range("A1").Select()
